# Regenerate the "K" (strikeouts) column (column G) of save_data sheet
# using the newly computed values (previously derived from "Strike#",
# now derived from K directly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..41 (row 1 is the header row)
$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 4
    8  = 2
    9  = 3
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 2
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 3
    25 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 2
    31 = 2
    32 = 3
    33 = 0
    34 = 2
    35 = 0
    36 = 1
    37 = 1
    38 = 3
    39 = 0
    40 = 1
    41 = 1
}

foreach ($row in ($newK.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
